$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.850.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +14.10%  "
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0691"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.120.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.858.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("E15").Value = "  +6.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.658"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.803.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0788"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.49%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("E34").Value = "  +5.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "91.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.35%  "
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.341.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.66%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("E42").Value = "  +8.55%  "
$ws.Range("E43").Value = "  +7.58%  "
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0519"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.017.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.45%  "
